$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 215, pushing the existing rows 215-220 down to 216-221.
# Excel will carry over the formatting (including the date style on column D)
# from the row above, matching the rest of the table.
$ws.Rows("215:215").Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(215, 1).Value = 5
$ws.Cells.Item(215, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(215, 3).Value = "Maule"
$ws.Cells.Item(215, 4).Value = 44509
$ws.Cells.Item(215, 5).Value = 7
$ws.Cells.Item(215, 6).Value = 100112023
$ws.Cells.Item(215, 7).Value = "Brócoli"
$ws.Cells.Item(215, 8).Value = "Sin especificar"
$ws.Cells.Item(215, 9).Value = "Primera"
$ws.Cells.Item(215, 10).Value = 5000
$ws.Cells.Item(215, 11).Value = 600
$ws.Cells.Item(215, 12).Value = 600
$ws.Cells.Item(215, 13).Value = 600
$ws.Cells.Item(215, 14).Value = "`$/unidad"
$ws.Cells.Item(215, 15).Value = "Región del Maule"
$ws.Cells.Item(215, 16).Value = 600
$ws.Cells.Item(215, 17).Value = 1
$ws.Cells.Item(215, 18).Value = "Hortaliza"
